# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# This script edits several worksheets of the listCalendar.211207 workbook:
#  - "Body"   (request body):  collapse the 3-field inline body into a single
#              schema reference row "listCalendar.211207Request".
#  - "200"    (response body): collapse the multi-field inline content into a
#              single schema reference row "listCalendar.211207Response".
#  - "400"    (response body): collapse the multi-field inline content into a
#              single schema reference row "errorResponse".
#  - "204"    gets a brand-new schema-reference row pointing at
#              "listCalendar.211207Response".
#  - "401", "403", "429", "500" each get a brand-new schema-reference row
#              pointing at the shared "errorResponse1" schema.
#  - "404" is left untouched.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param($ws, $Row, $Section, $Name, $SchemaName, $Mandatory)

    $ws.Cells.Item($Row, 1).Value  = $Section      # A -> Section
    $ws.Cells.Item($Row, 2).Value  = $Name          # B -> Name
    $ws.Cells.Item($Row, 3).ClearContents()         # C -> Parent
    $ws.Cells.Item($Row, 4).ClearContents()         # D -> Description
    $ws.Cells.Item($Row, 5).Value  = "schema"       # E -> Type
    $ws.Cells.Item($Row, 6).ClearContents()         # F -> Items Data Type
    $ws.Cells.Item($Row, 7).Value  = $SchemaName    # G -> Schema Name
    $ws.Cells.Item($Row, 8).ClearContents()         # H -> Format
    $ws.Cells.Item($Row, 9).Value  = $Mandatory     # I -> Mandatory
    $ws.Cells.Item($Row, 10).ClearContents()        # J -> Min Value/Length/Item
    $ws.Cells.Item($Row, 11).ClearContents()        # K -> Max Value/Length/Item
    $ws.Cells.Item($Row, 12).ClearContents()        # L -> PatternEba
    $ws.Cells.Item($Row, 13).ClearContents()        # M -> Regex
    $ws.Cells.Item($Row, 14).ClearContents()        # N -> Allowed value
    $ws.Cells.Item($Row, 15).ClearContents()        # O -> Example
}

# ---------------------------------------------------------------------------
# "Body" sheet: row 3 becomes the schema reference, rows 4-5 (settlementBIC,
# businessDate) go away entirely.
# ---------------------------------------------------------------------------
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow $wsBody 3 "body" "listCalendar.211207Request" "listCalendar.211207Request" "Yes"
$wsBody.Rows("4:5").Delete()

# ---------------------------------------------------------------------------
# "200" sheet: row 3 becomes the schema reference, rows 4-11 (businessDate,
# exceptionLacValues and its nested fields) go away entirely.
# ---------------------------------------------------------------------------
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow $ws200 3 "content" "listCalendar.211207Response" "listCalendar.211207Response" "Yes"
$ws200.Rows("4:11").Delete()

# ---------------------------------------------------------------------------
# "204" sheet: brand-new row 3 referencing the same response schema.
# ---------------------------------------------------------------------------
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow $ws204 3 "content" "listCalendar.211207Response" "listCalendar.211207Response" "Yes"

# ---------------------------------------------------------------------------
# "400" sheet: row 3 becomes the schema reference, rows 4-6 (errorCode,
# errorCodeDescription, requestId) go away entirely.
# ---------------------------------------------------------------------------
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow $ws400 3 "content" "errorResponse" "errorResponse" "Yes"
$ws400.Rows("4:6").Delete()

# ---------------------------------------------------------------------------
# "401", "403", "429", "500": brand-new row 3 referencing the shared
# "errorResponse1" schema. "404" is intentionally left unchanged.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("401", "403", "429", "500")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Set-SchemaRow $ws 3 "content" "errorResponse1" "errorResponse1" "Yes"
}
